# Refresh market-price derived columns (H:N) across the Leve profit sheets.
# Values mirror an external market-data snapshot (currentAveragePrice* / Leve*Price* / Leve*Profit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 989.8570999999999
$ws.Range("I32").Value = 999
$ws.Range("J32").Value = 988.3333
$ws.Range("K32").Value = 999
$ws.Range("L32").Value = 988.3333
$ws.Range("M32").Value = -673
$ws.Range("N32").Value = -1640.3333

$ws.Range("H98").Value = 929.4706
$ws.Range("I98").Value = 816.73334
$ws.Range("J98").Value = 1775
$ws.Range("K98").Value = 816.73334
$ws.Range("L98").Value = 1775
$ws.Range("M98").Value = 681.26666
$ws.Range("N98").Value = -4771

$ws.Range("H122").Value = 929.4706
$ws.Range("I122").Value = 816.73334
$ws.Range("J122").Value = 1775
$ws.Range("K122").Value = 2450.20002
$ws.Range("L122").Value = 5325
$ws.Range("M122").Value = -0.2000200000002224
$ws.Range("N122").Value = -10225

$ws.Range("H129").Value = 1068.24
$ws.Range("I129").Value = 567.6923
$ws.Range("J129").Value = 1244.1082
$ws.Range("K129").Value = 1703.0769
$ws.Range("L129").Value = 3732.3246
$ws.Range("M129").Value = 3296.9231
$ws.Range("N129").Value = -13732.3246

$ws.Range("H138").Value = 2533.756
$ws.Range("I138").Value = 2215.8948
$ws.Range("J138").Value = 2808.2727
$ws.Range("K138").Value = 6647.6844
$ws.Range("L138").Value = 8424.8181
$ws.Range("M138").Value = -1507.6844
$ws.Range("N138").Value = -18704.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 450007.16
$ws.Range("I32").Value = 546920.2
$ws.Range("J32").Value = 13898.667
$ws.Range("K32").Value = 546920.2
$ws.Range("L32").Value = 13898.667
$ws.Range("M32").Value = -546633.2
$ws.Range("N32").Value = -14472.667

$ws.Range("H45").Value = 3546.5715
$ws.Range("I45").Value = 2618.6667
$ws.Range("J45").Value = 4242.5
$ws.Range("K45").Value = 2618.6667
$ws.Range("L45").Value = 4242.5
$ws.Range("M45").Value = -2241.6667
$ws.Range("N45").Value = -4996.5

$ws.Range("H107").Value = 28150
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 28150
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 28150
$ws.Range("N107").Value = -35830

$ws.Range("H122").Value = 1449.0588
$ws.Range("I122").Value = 1292.7693
$ws.Range("J122").Value = 1957
$ws.Range("K122").Value = 3878.3079
$ws.Range("L122").Value = 5871
$ws.Range("M122").Value = -1428.3079
$ws.Range("N122").Value = -10771

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2208.423
$ws.Range("I20").Value = 2030.1538
$ws.Range("J20").Value = 2386.6924
$ws.Range("K20").Value = 2030.1538
$ws.Range("L20").Value = 2386.6924
$ws.Range("M20").Value = -1783.1538
$ws.Range("N20").Value = -2880.6924

$ws.Range("H80").Value = 1078.762
$ws.Range("I80").Value = 2374.889
$ws.Range("J80").Value = 106.666664
$ws.Range("K80").Value = 2374.889
$ws.Range("L80").Value = 106.666664
$ws.Range("M80").Value = -1376.889
$ws.Range("N80").Value = -2102.666664

$ws.Range("H83").Value = 1078.762
$ws.Range("I83").Value = 2374.889
$ws.Range("J83").Value = 106.666664
$ws.Range("K83").Value = 11874.445
$ws.Range("L83").Value = 533.33332
$ws.Range("M83").Value = -6882.445
$ws.Range("N83").Value = -10517.33332

$ws.Range("H140").Value = 97106.664
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 97106.664
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 97106.664
$ws.Range("N140").Value = -107466.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H122").Value = 1983.8422
$ws.Range("I122").Value = 1897.6666
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5692.9998
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3242.9998
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 4388104.5
$ws.Range("I132").Value = 1655.25
$ws.Range("J132").Value = 9261937
$ws.Range("K132").Value = 4965.75
$ws.Range("L132").Value = 27785811
$ws.Range("M132").Value = -2435.75
$ws.Range("N132").Value = -27790871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1269.03
$ws.Range("I68").Value = 661.6799999999999
$ws.Range("J68").Value = 1471.48
$ws.Range("K68").Value = 1985.04
$ws.Range("L68").Value = 4414.440000000001
$ws.Range("M68").Value = -1174.04
$ws.Range("N68").Value = -6036.440000000001

$ws.Range("H71").Value = 1269.03
$ws.Range("I71").Value = 661.6799999999999
$ws.Range("J71").Value = 1471.48
$ws.Range("K71").Value = 5955.12
$ws.Range("L71").Value = 13243.32
$ws.Range("M71").Value = -1899.12
$ws.Range("N71").Value = -21355.32

$ws.Range("H98").Value = 500250
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 500250
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 1500750
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -1503746

$ws.Range("H107").Value = 1945.175
$ws.Range("I107").Value = 250.6875
$ws.Range("J107").Value = 3074.8333
$ws.Range("K107").Value = 752.0625
$ws.Range("L107").Value = 9224.499899999999
$ws.Range("M107").Value = 1167.9375
$ws.Range("N107").Value = -13064.4999

$ws.Range("H132").Value = 3462.5715
$ws.Range("I132").Value = 2553
$ws.Range("J132").Value = 4101.7295
$ws.Range("K132").Value = 22977
$ws.Range("L132").Value = 36915.5655
$ws.Range("M132").Value = -20447
$ws.Range("N132").Value = -41975.5655

$ws.Range("H137").Value = 9003.058999999999
$ws.Range("I137").Value = 14428.777
$ws.Range("J137").Value = 2899.125
$ws.Range("K137").Value = 43286.331
$ws.Range("L137").Value = 8697.375
$ws.Range("M137").Value = -38186.331
$ws.Range("N137").Value = -18897.375

$ws.Range("H139").Value = 2408.5715
$ws.Range("I139").Value = 2408.5715
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7225.7145
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2085.7145

$ws.Range("H140").Value = 1783.5
$ws.Range("I140").Value = 1223.8
$ws.Range("J140").Value = 2902.9
$ws.Range("K140").Value = 3671.4
$ws.Range("L140").Value = 8708.700000000001
$ws.Range("M140").Value = 1508.6
$ws.Range("N140").Value = -19068.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5398.8335
$ws.Range("I122").Value = 7003
$ws.Range("J122").Value = 5078
$ws.Range("K122").Value = 21009
$ws.Range("L122").Value = 15234
$ws.Range("M122").Value = -18559
$ws.Range("N122").Value = -20134

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 33337096
$ws.Range("I40").Value = 62502172
$ws.Range("J40").Value = 5578.5713
$ws.Range("K40").Value = 62502172
$ws.Range("L40").Value = 5578.5713
$ws.Range("M40").Value = -62502036
$ws.Range("N40").Value = -5850.5713

$ws.Range("H122").Value = 6100.727
$ws.Range("I122").Value = 4027
$ws.Range("J122").Value = 7285.7144
$ws.Range("K122").Value = 12081
$ws.Range("L122").Value = 21857.1432
$ws.Range("M122").Value = -9631
$ws.Range("N122").Value = -26757.1432

$ws.Range("H132").Value = 4138.9443
$ws.Range("I132").Value = 3855.9583
$ws.Range("J132").Value = 4704.9165
$ws.Range("K132").Value = 11567.8749
$ws.Range("L132").Value = 14114.7495
$ws.Range("M132").Value = -9037.874899999999
$ws.Range("N132").Value = -19174.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 1960
$ws.Range("I122").Value = 1572.7273
$ws.Range("J122").Value = 3025
$ws.Range("K122").Value = 4718.1819
$ws.Range("L122").Value = 9075
$ws.Range("M122").Value = -2268.1819
$ws.Range("N122").Value = -13975

$ws.Range("H126").Value = 1439.1333
$ws.Range("I126").Value = 1216.3334
$ws.Range("J126").Value = 1587.6666
$ws.Range("K126").Value = 3649.0002
$ws.Range("L126").Value = 4762.9998
$ws.Range("M126").Value = -1179.0002
$ws.Range("N126").Value = -9702.9998

$ws.Range("H132").Value = 3032049
$ws.Range("I132").Value = 1486
$ws.Range("J132").Value = 8774168
$ws.Range("K132").Value = 4458
$ws.Range("L132").Value = 26322504
$ws.Range("M132").Value = -1928
$ws.Range("N132").Value = -26327564
